$wb = $excel.ActiveWorkbook

# The file "163d6d51-a136-4906-a1a5-d29efb81317b.md" (row 4 in every sheet) is
# "Ready for handoff". Generating the handoff report stamps a fresh handoff
# datetime for that row on the Overview sheet as well as on each of the
# per-language (zh-cn / de-de) sheets.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value = "2016-58-14 08:58:04"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-14 08:57:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-14 08:58:04"
